# Add 2022-Q4 data:
#  - insert a new "2022-Q4" sheet right after "总计" (pushing "2021-Q4" to 3rd position)
#  - populate "2022-Q4" with its fund-holding header/data row
#  - insert a new summary row on "总计" for 2022-Q4, shifting the old 2021-Q4
#    summary row down to row 3
#
# NOTE: worksheet object references captured *before* a sheet-collection
# mutation (Add/Move/...) can go stale and silently resolve to whatever
# sheet now occupies their old slot, so every sheet reference below is
# (re-)fetched by name *after* the mutating calls that could invalidate it.

$wb = $excel.ActiveWorkbook

# ---- create + position the new "2022-Q4" sheet --------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"

$total = $wb.Worksheets.Item("总计")
$newSheet.Move($null, $total)

# re-fetch all sheet handles fresh, now that the sheet collection is stable
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item("2022-Q4")

# ---- populate "2022-Q4" (fund holdings for the quarter) -----------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# header row + index column share the same bold / thin-border / centered look
# used throughout the rest of the workbook (comma multi-area Range doesn't
# reliably style every area here, so style each piece separately)
foreach ($rng in @($newSheet.Range("B1:H1"), $newSheet.Range("A2"))) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "000822"
$newSheet.Range("C2").Value = "东海美丽中国灵活配置混合"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.12"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "76.79"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "2.14"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0026"

$newSheet.Range("H2").Value = 8

# ---- update the "总计" summary sheet -------------------------------------
# shift the existing 2021-Q4 summary row from row 2 down to row 3
# (match A2's existing bold / thin-border / centered index-column look)
$a3 = $total.Range("A3")
$a3.Font.Bold = $true
$a3.Borders.LineStyle = 1
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4160
$a3.Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 2.88

# write the new 2022-Q4 summary row in row 2
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0
